# Georgia_Converted.xlsx edit script
# - Turns off the "travel_limit"(U7) weight (1 -> 0) and updates the
#   LockdownEffectiveness denominator (X7) from 13 to 12.
# - Recomputes the weighted "LockdownEffectiveness" column X for every
#   existing date row (26-221) using the updated weights/denominator.
# - Appends 12 new date rows (222-233) for 9/30/2020 .. 10/11/2020,
#   copying the policy pattern of the last existing row (221) and
#   recomputing their X values the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1. Update the weight row (row 7): travel_limit (U7) goes from 1 to 0,
#    and the LockdownEffectiveness denominator (X7) goes from 13 to 12.
# ---------------------------------------------------------------
$ws.Range("U7").Value = 0
$ws.Range("X7").Value = 12

# Read back the full updated weight vector B7:W7 (22 values).
$weights = @()
for ($c = 2; $c -le 23; $c++) {
    $weights += $ws.Cells.Item(7, $c).Value2
}
$denom = $ws.Range("X7").Value2

# ---------------------------------------------------------------
# 2. Recompute column X (LockdownEffectiveness) for every existing
#    date row (26 through 221) using the new weights/denominator.
# ---------------------------------------------------------------
for ($r = 26; $r -le 221; $r++) {
    $sum = 0.0
    for ($c = 2; $c -le 23; $c++) {
        $sum += $weights[$c - 2] * $ws.Cells.Item($r, $c).Value2
    }
    $ws.Cells.Item($r, 24).Value = $sum / $denom
}

# ---------------------------------------------------------------
# 3. Append 12 new date rows (222-233) for 9/30/2020 .. 10/11/2020.
#    They repeat the same policy pattern as the last existing row (221).
# ---------------------------------------------------------------
$newDates = @("9/30/2020","10/1/2020","10/2/2020","10/3/2020","10/4/2020","10/5/2020","10/6/2020","10/7/2020","10/8/2020","10/9/2020","10/10/2020","10/11/2020")
$firstNewRow = 222
$lastNewRow = $firstNewRow + $newDates.Count - 1

# Copy formatting (font/border/alignment/number-format) of the last
# existing data row down across all the new rows in one shot.
$ws.Range("A221:X221").Copy()
$ws.Range("A" + $firstNewRow + ":X" + $lastNewRow).PasteSpecial(-4122)

# Write each date as a literal text formula (="9/30/2020") so Excel
# doesn't auto-convert the text to a date serial, then flatten the
# formulas down to plain shared-string values.
for ($i = 0; $i -lt $newDates.Count; $i++) {
    $r = $firstNewRow + $i
    $ws.Cells.Item($r, 1).Formula = "=""" + $newDates[$i] + """"
}
$ws.Range("A" + $firstNewRow + ":A" + $lastNewRow).Copy()
$ws.Range("A" + $firstNewRow + ":A" + $lastNewRow).PasteSpecial(-4163)

# Same policy (B:W) pattern as row 221, replicated on every new row,
# followed by the recomputed X value.
$rowVals = @(0,0,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,0)
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    for ($c = 2; $c -le 23; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 2]
    }
    $sum = 0.0
    for ($c = 2; $c -le 23; $c++) {
        $sum += $weights[$c - 2] * $rowVals[$c - 2]
    }
    $ws.Cells.Item($r, 24).Value = $sum / $denom
}
